# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their text values (e.g. "215.39", "7.55")
# instead of being auto-converted to numbers by Excel's smart parsing.
$priceCells = "D2","D3","D5","D8","D10","D11","D12","D13","D15","D16","D17","D18","D20","D25","D26","D27","D28","D29","D33","D36","D38","D40","D44","D45","D46","D47","D48","D49","D51"
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.613.67"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.666.31"
$ws.Range("E3").Value = "  -3.49%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "215.39"
$ws.Range("E5").Value = "  -1.83%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Solana
$ws.Range("D8").Value = "23.65"
$ws.Range("E8").Value = "  -2.22%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.78%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0622"
$ws.Range("E10").Value = "  -2.11%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  -2.12%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.901.51"
$ws.Range("E12").Value = "  -3.52%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.671.05"
$ws.Range("E13").Value = "  -3.14%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.47%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.558"
$ws.Range("E15").Value = "  -1.48%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "66.23"
$ws.Range("E16").Value = "  -2.47%  "

# Row 17 - BitcoinCash
$ws.Range("D17").Value = "247.67"
$ws.Range("E17").Value = "  +1.75%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "27.627.66"
$ws.Range("E18").Value = "  -1.13%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -3.66%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -4.65%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.66%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -5.16%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -5.00%  "

# Row 25 - Monero
$ws.Range("D25").Value = "146.01"
$ws.Range("E25").Value = "  -2.15%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "7.18"
$ws.Range("E26").Value = "  -4.81%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "16.44"
$ws.Range("E27").Value = "  -2.33%  "

# Row 28 and 29 swap places: BinanceUSD (was 28) <-> Stellar (was 29)
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +3.55%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.96%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -3.06%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.478.11"
$ws.Range("E33").Value = "  -0.65%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -5.43%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -6.06%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "0.939"
$ws.Range("E36").Value = "  -2.54%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -1.21%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.574"
$ws.Range("E38").Value = "  -6.15%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.44%  "

# Row 40 - Aave
$ws.Range("D40").Value = "69.50"
$ws.Range("E40").Value = "  -2.74%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -5.55%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.05%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -7.34%  "

# Row 44 and 45 swap places: RocketPoolETH (was 44) <-> MXToken (was 45)
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.21"
$ws.Range("E44").Value = "  -3.86%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.809.12"
$ws.Range("E45").Value = "  -3.48%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = "0.787"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.71"
$ws.Range("E47").Value = "  -3.63%  "

# Row 48 - Quant
$ws.Range("D48").Value = "89.27"
$ws.Range("E48").Value = "  -2.86%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -3.66%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -3.11%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.90"
$ws.Range("E51").Value = "  -3.98%  "
